# Auto-generated Excel COM-interop script
# Updates market-price-derived columns (H-N) on multiple worksheets
# to match refreshed pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 790.6667
$ws.Range("I11").Value = 790.6667
$ws.Range("K11").Value = 790.6667
$ws.Range("M11").Value = -650.6667

$ws.Range("H43").Value = 515561.75
$ws.Range("J43").Value = 824599.8
$ws.Range("L43").Value = 824599.8
$ws.Range("N43").Value = -824737.8

$ws.Range("H58").Value = 35714700
$ws.Range("J58").Value = 83334070
$ws.Range("L58").Value = 250002210
$ws.Range("N58").Value = -250002510

$ws.Range("H98").Value = 3175.762
$ws.Range("I98").Value = 3175.762
$ws.Range("K98").Value = 3175.762
$ws.Range("M98").Value = -1677.762

$ws.Range("H116").Value = 14716679
$ws.Range("I116").Value = 35721790
$ws.Range("K116").Value = 35721790
$ws.Range("M116").Value = -35718348

$ws.Range("H122").Value = 3175.762
$ws.Range("I122").Value = 3175.762
$ws.Range("K122").Value = 9527.286
$ws.Range("M122").Value = -7077.286

$ws.Range("H123").Value = 52465.4
$ws.Range("J123").Value = 52465.4
$ws.Range("L123").Value = 52465.4
$ws.Range("N123").Value = -62265.4

$ws.Range("H132").Value = 1083.7142
$ws.Range("I132").Value = 1027.9
$ws.Range("K132").Value = 3083.7
$ws.Range("M132").Value = -553.7000000000003

$ws.Range("H138").Value = 1566225.4
$ws.Range("I138").Value = 1695.6072
$ws.Range("K138").Value = 5086.821599999999
$ws.Range("M138").Value = 53.17840000000069

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4088349.2
$ws.Range("J32").Value = 27249.5
$ws.Range("L32").Value = 27249.5
$ws.Range("N32").Value = -27823.5

$ws.Range("H34").Value = 94041.25
$ws.Range("J34").Value = 59249.5
$ws.Range("L34").Value = 59249.5
$ws.Range("N34").Value = -59791.5

$ws.Range("H45").Value = 2053.9524
$ws.Range("I45").Value = 2260.4
$ws.Range("J45").Value = 1866.2727
$ws.Range("K45").Value = 2260.4
$ws.Range("L45").Value = 1866.2727
$ws.Range("M45").Value = -1883.4
$ws.Range("N45").Value = -2620.2727

$ws.Range("H61").Value = 45457772
$ws.Range("I61").Value = 1504.2307
$ws.Range("J61").Value = 111116824
$ws.Range("K61").Value = 1504.2307
$ws.Range("L61").Value = 111116824
$ws.Range("M61").Value = -1292.2307
$ws.Range("N61").Value = -111117248

$ws.Range("H97").Value = 2382320.5
$ws.Range("I97").Value = 1337.36
$ws.Range("J97").Value = 8334778.5
$ws.Range("K97").Value = 1337.36
$ws.Range("L97").Value = 8334778.5
$ws.Range("M97").Value = -841.3599999999999
$ws.Range("N97").Value = -8335770.5

$ws.Range("H132").Value = 4359.0864
$ws.Range("I132").Value = 3555.475
$ws.Range("J132").Value = 6144.8887
$ws.Range("K132").Value = 10666.425
$ws.Range("L132").Value = 18434.6661
$ws.Range("M132").Value = -8136.424999999999
$ws.Range("N132").Value = -23494.6661

$ws.Range("H136").Value = 45457772
$ws.Range("I136").Value = 1504.2307
$ws.Range("J136").Value = 111116824
$ws.Range("K136").Value = 4512.6921
$ws.Range("L136").Value = 333350472
$ws.Range("M136").Value = -1962.6921
$ws.Range("N136").Value = -333355572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1698.6666
$ws.Range("I94").Value = 718.4375
$ws.Range("K94").Value = 718.4375
$ws.Range("M94").Value = -267.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 71428740
$ws.Range("I7").Value = 89.5
$ws.Range("J7").Value = 166666940
$ws.Range("K7").Value = 89.5
$ws.Range("L7").Value = 166666940
$ws.Range("M7").Value = 23.5
$ws.Range("N7").Value = -166667166

$ws.Range("H59").Value = 58997.332
$ws.Range("J59").Value = 58997.332
$ws.Range("L59").Value = 58997.332
$ws.Range("N59").Value = -61287.332

$ws.Range("H122").Value = 4687.4116
$ws.Range("I122").Value = 4142.943
$ws.Range("K122").Value = 12428.829
$ws.Range("M122").Value = -9978.829000000002

$ws.Range("H132").Value = 3520.878
$ws.Range("I132").Value = 2692.2222
$ws.Range("K132").Value = 8076.6666
$ws.Range("M132").Value = -5546.6666

$ws.Range("H134").Value = 3878.0557
$ws.Range("I134").Value = 2442.7827
$ws.Range("J134").Value = 6417.385
$ws.Range("K134").Value = 7328.348100000001
$ws.Range("L134").Value = 19252.155
$ws.Range("M134").Value = -4793.348100000001
$ws.Range("N134").Value = -24322.155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 857.5
$ws.Range("I5").Value = 574.8570999999999
$ws.Range("J5").Value = 1846.75
$ws.Range("K5").Value = 1724.5713
$ws.Range("L5").Value = 5540.25
$ws.Range("M5").Value = -1612.5713
$ws.Range("N5").Value = -5764.25

$ws.Range("H11").Value = 2724.75
$ws.Range("I11").Value = 1799.6666
$ws.Range("K11").Value = 5398.9998
$ws.Range("M11").Value = -5258.9998

$ws.Range("H113").Value = 3924.2778
$ws.Range("J113").Value = 5220.231
$ws.Range("L113").Value = 15660.693
$ws.Range("N113").Value = -20000.693

$ws.Range("H135").Value = 857.5
$ws.Range("I135").Value = 574.8570999999999
$ws.Range("J135").Value = 1846.75
$ws.Range("K135").Value = 5173.7139
$ws.Range("L135").Value = 16620.75
$ws.Range("M135").Value = -2638.7139
$ws.Range("N135").Value = -21690.75

$ws.Range("H137").Value = 106848.42
$ws.Range("I137").Value = 77805.766
$ws.Range("K137").Value = 233417.298
$ws.Range("M137").Value = -228317.298

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 20000
$ws.Range("J26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("N26").Value = -20560

$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -20996

$ws.Range("H52").Value = 83749.75
$ws.Range("J52").Value = 90000
$ws.Range("L52").Value = 90000
$ws.Range("N52").Value = -90518

$ws.Range("H58").Value = 57277.285
$ws.Range("J58").Value = 78980
$ws.Range("L58").Value = 78980
$ws.Range("N58").Value = -79534

$ws.Range("H63").Value = 54900
$ws.Range("J63").Value = 54900
$ws.Range("L63").Value = 54900
$ws.Range("N63").Value = -56272

$ws.Range("H64").Value = 69725.2
$ws.Range("J64").Value = 99313
$ws.Range("L64").Value = 99313
$ws.Range("N64").Value = -99809

$ws.Range("H66").Value = 54900
$ws.Range("J66").Value = 54900
$ws.Range("L66").Value = 164700
$ws.Range("N66").Value = -171564

$ws.Range("H67").Value = 69725.2
$ws.Range("J67").Value = 99313
$ws.Range("L67").Value = 99313
$ws.Range("N67").Value = -101029

$ws.Range("H70").Value = 7291.7617
$ws.Range("I70").Value = 4818.5
$ws.Range("K70").Value = 4818.5
$ws.Range("M70").Value = -4548.5

$ws.Range("H73").Value = 7291.7617
$ws.Range("I73").Value = 4818.5
$ws.Range("K73").Value = 4818.5
$ws.Range("M73").Value = -3882.5

$ws.Range("H75").Value = 28994.777
$ws.Range("J75").Value = 28994.777
$ws.Range("L75").Value = 28994.777
$ws.Range("N75").Value = -30742.777

$ws.Range("H78").Value = 28994.777
$ws.Range("J78").Value = 28994.777
$ws.Range("L78").Value = 86984.33099999999
$ws.Range("N78").Value = -95720.33099999999

$ws.Range("H97").Value = 998.15
$ws.Range("I97").Value = 1132
$ws.Range("J97").Value = 797.375
$ws.Range("K97").Value = 1132
$ws.Range("L97").Value = 797.375
$ws.Range("M97").Value = -636
$ws.Range("N97").Value = -1789.375

$ws.Range("H102").Value = 3160.9722
$ws.Range("I102").Value = 3020.0625
$ws.Range("J102").Value = 4288.25
$ws.Range("K102").Value = 3020.0625
$ws.Range("L102").Value = 4288.25
$ws.Range("M102").Value = -1398.0625
$ws.Range("N102").Value = -7532.25

$ws.Range("H113").Value = 6004.614
$ws.Range("I113").Value = 2255.0476
$ws.Range("K113").Value = 2255.0476
$ws.Range("M113").Value = -85.04759999999987

$ws.Range("H122").Value = 1813582.4
$ws.Range("I122").Value = 2339170.5
$ws.Range("K122").Value = 7017511.5
$ws.Range("M122").Value = -7015061.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 16678332
$ws.Range("J33").Value = 17500
$ws.Range("L33").Value = 17500
$ws.Range("N33").Value = -18080

$ws.Range("H122").Value = 4360.7827
$ws.Range("I122").Value = 3504.7273
$ws.Range("K122").Value = 10514.1819
$ws.Range("M122").Value = -8064.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 90999.336
$ws.Range("J46").Value = 90999.336
$ws.Range("L46").Value = 90999.336
$ws.Range("N46").Value = -91461.336

$ws.Range("H132").Value = 4355.1284
$ws.Range("I132").Value = 4172.4194
$ws.Range("K132").Value = 12517.2582
$ws.Range("M132").Value = -9987.2582

$ws.Range("H134").Value = 90999.336
$ws.Range("J134").Value = 90999.336
$ws.Range("L134").Value = 272998.008
$ws.Range("N134").Value = -278068.008

$ws.Range("H136").Value = 22957874
$ws.Range("I136").Value = 45455348
$ws.Range("K136").Value = 136366044
$ws.Range("M136").Value = -136363494
